# Update newsbot state: append rows 96-99 to the "Historico" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96
$ws.Range("A96").Value = "05/01/2026 14:45:58"
$ws.Range("B96").Value = "05/01 14:22"
$ws.Range("C96").Value = "g1 > Política"
$ws.Range("D96").Value = "Brasil condena intervenção armada na Venezuela:  'Não podemos aceitar o argumento de que os fins justificam os meios'"
$ws.Range("E96").Value = "https://g1.globo.com/politica/noticia/2026/01/05/brasil-discursa-no-conselho-de-seguranca-da-onu-nao-podemos-aceitar-o-argumento-de-que-os-fins-justificam-os-meios.ghtml"
$ws.Range("F96").Value = "ldo"
$ws.Range("G96").Value = "bilidade de outorgar aos mais fortes o que é justo, injusto, o que é correto.`nO blog do Va&lt;b&gt;ldo&lt;/b&gt; Cruz já tinha adiantado a informação sobre a posição do Brasil. `nVeja os vídeos que estão"

# Row 97
$ws.Range("A97").Value = "05/01/2026 14:45:59"
$ws.Range("B97").Value = "05/01 14:17"
$ws.Range("C97").Value = "Metrópoles"
$ws.Range("D97").Value = "Ministro do TCU diz que nota técnica do BC sobre liquidação do Master carece de prova documental"
$ws.Range("E97").Value = "https://www.metropoles.com/colunas/tacio-lorran/caso-master-ministro-do-tcu-aponta-falta-de-documentos-em-nota-do-bc"
$ws.Range("F97").Value = "tcu"
$ws.Range("G97").Value = "Jhonatan de Jesus determinou que área técnica do TCU inspecione documentos em posse do BC sobre liquidação do Master"

# Row 98 (no value published-time in column B)
$ws.Range("A98").Value = "05/01/2026 14:45:59"
$ws.Range("C98").Value = "VEJA"
$ws.Range("D98").Value = "Bancos defendem BC em conflito com TCU no caso Master"
$ws.Range("E98").Value = "https://veja.abril.com.br/economia/bancos-defendem-bc-em-conflito-com-tcu-no-caso-master/"
$ws.Range("F98").Value = "banco central"
$ws.Range("G98").Value = "Tribunal de Contas da União questinou decisão do Banco Central que fechou banco de Daniel Vorcaro"

# Row 99 (B, F, G are present but empty)
$ws.Range("A99").Value = "05/01/2026 14:46:00"
$ws.Range("B99").Value = ""
$ws.Range("C99").Value = "VEJA"
$ws.Range("D99").Value = "Vale à pena parcelar o IPVA? Entenda qual a melhor forma de pagar o imposto"
$ws.Range("E99").Value = "https://veja.abril.com.br/economia/vale-a-pena-parcelar-o-ipva-entenda-qual-a-melhor-forma-de-pagar-o-imposto/"
$ws.Range("F99").Value = ""
$ws.Range("G99").Value = ""
